$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-20 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-21 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("48×29=", $true, $false, $false, $false, $false, $true, 1, $false, "91×19=", 2) | Out-Null
$d.Content.Find.Execute("19×68=", $true, $false, $false, $false, $false, $true, 1, $false, "96×73=", 2) | Out-Null
$d.Content.Find.Execute("12×51=", $true, $false, $false, $false, $false, $true, 1, $false, "38×28=", 2) | Out-Null
$d.Content.Find.Execute("60×98=", $true, $false, $false, $false, $false, $true, 1, $false, "24×13=", 2) | Out-Null
$d.Content.Find.Execute("27×66=", $true, $false, $false, $false, $false, $true, 1, $false, "72×47=", 2) | Out-Null
$d.Content.Find.Execute("46×95=", $true, $false, $false, $false, $false, $true, 1, $false, "70×45=", 2) | Out-Null
$d.Content.Find.Execute("85×13=", $true, $false, $false, $false, $false, $true, 1, $false, "34×80=", 2) | Out-Null
$d.Content.Find.Execute("45×96=", $true, $false, $false, $false, $false, $true, 1, $false, "87×84=", 2) | Out-Null
$d.Content.Find.Execute("47×35=", $true, $false, $false, $false, $false, $true, 1, $false, "96×69=", 2) | Out-Null
$d.Content.Find.Execute("14×90=", $true, $false, $false, $false, $false, $true, 1, $false, "20×58=", 2) | Out-Null
$d.Content.Find.Execute("68×64=", $true, $false, $false, $false, $false, $true, 1, $false, "76×84=", 2) | Out-Null
$d.Content.Find.Execute("26×83=", $true, $false, $false, $false, $false, $true, 1, $false, "45×63=", 2) | Out-Null
$d.Content.Find.Execute("41×65=", $true, $false, $false, $false, $false, $true, 1, $false, "96×58=", 2) | Out-Null
$d.Content.Find.Execute("52×99=", $true, $false, $false, $false, $false, $true, 1, $false, "64×85=", 2) | Out-Null
$d.Content.Find.Execute("34×85=", $true, $false, $false, $false, $false, $true, 1, $false, "34×63=", 2) | Out-Null
$d.Content.Find.Execute("98×37=", $true, $false, $false, $false, $false, $true, 1, $false, "57×54=", 2) | Out-Null
$d.Content.Find.Execute("46×54=", $true, $false, $false, $false, $false, $true, 1, $false, "84×22=", 2) | Out-Null
$d.Content.Find.Execute("64×34=", $true, $false, $false, $false, $false, $true, 1, $false, "37×22=", 2) | Out-Null
$d.Content.Find.Execute("20×32=", $true, $false, $false, $false, $false, $true, 1, $false, "85×93=", 2) | Out-Null
$d.Content.Find.Execute("45×62=", $true, $false, $false, $false, $false, $true, 1, $false, "80×69=", 2) | Out-Null
$d.Content.Find.Execute("14×39=", $true, $false, $false, $false, $false, $true, 1, $false, "91×72=", 2) | Out-Null
$d.Content.Find.Execute("81×54=", $true, $false, $false, $false, $false, $true, 1, $false, "72×50=", 2) | Out-Null
$d.Content.Find.Execute("46×15=", $true, $false, $false, $false, $false, $true, 1, $false, "45×28=", 2) | Out-Null
$d.Content.Find.Execute("88×47=", $true, $false, $false, $false, $false, $true, 1, $false, "34×73=", 2) | Out-Null
$d.Content.Find.Execute("78×88=", $true, $false, $false, $false, $false, $true, 1, $false, "86×73=", 2) | Out-Null
